$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.988.54"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "2.587.30"
$ws.Range("E3").Value = "  +1.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'583.80"
$ws.Range("E5").Value = "  +1.55%  "

# Row 6
$ws.Range("D6").Value = "'147.49"
$ws.Range("E6").Value = "  +0.81%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +2.73%  "

# Row 9
$ws.Range("E9").Value = "  +2.89%  "

# Row 10
$ws.Range("D10").Value = "'5.67"
$ws.Range("E10").Value = "  +3.19%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("E12").Value = "  -0.09%  "

# Row 13
$ws.Range("D13").Value = "'27.43"
$ws.Range("E13").Value = "  +1.33%  "

# Row 14
$ws.Range("D14").Value = "3.050.02"
$ws.Range("E14").Value = "  +1.69%  "

# Row 15
$ws.Range("D15").Value = "62.841.01"

# Row 16
$ws.Range("E16").Value = "  +3.65%  "

# Row 17
$ws.Range("D17").Value = "2.584.78"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18
$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").Value = "'342.26"
$ws.Range("E19").Value = "  +2.07%  "

# Row 20
$ws.Range("E20").Value = "  +2.13%  "

# Row 21
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -0.63%  "

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'67.21"
$ws.Range("E23").Value = "  +3.11%  "

# Row 24
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "2.713.16"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -1.55%  "

# Row 26
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'1.60"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'7.89"
$ws.Range("E28").Value = "  +8.85%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.32"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "'1.44"
$ws.Range("E30").Value = "  -1.86%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.93"
$ws.Range("E31").Value = "  +3.61%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0823"
$ws.Range("E32").Value = "  +1.74%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'462.87"
$ws.Range("E33").Value = "  +14.26%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'175.03"
$ws.Range("E34").Value = "  -1.49%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +4.40%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.08%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.403"
$ws.Range("E37").Value = "  +1.24%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.10"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'4.55"
$ws.Range("E39").Value = "  +5.08%  "

# Row 40
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = "  -1.26%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'159.25"
$ws.Range("E42").Value = "  +5.57%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "  +0.95%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.639"
$ws.Range("E44").Value = "  +6.51%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'21.25"
$ws.Range("E45").Value = "  +2.34%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0543"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0968"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0237"
$ws.Range("E48").Value = "  -0.45%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'18.42"
$ws.Range("E49").Value = "  +1.68%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  +1.15%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.42"
$ws.Range("E51").Value = "  +1.07%  "
